$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 1536869
$ws.Range("B10").Value = "Test One"
$ws.Range("C10").Value = "These guys are good"
